# "scale data to match units" — rescale the blood-data columns on the
# active ("data") sheet so everything lines up with the AB-column units,
# add a formatted (but empty) marker cell, and leave the selection where
# the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Y/Z/AA get rescaled (AB12 is already in the correct units and
# is left untouched).
$ws.Range("Y12").Value  = 2.6038257396855347
$ws.Range("Z12").Value  = 2.6038257396855347
$ws.Range("AA12").Value = 2.4202870927432429

# Rows 23-26: Z/AA get rescaled, and the now-redundant AB column (which
# duplicated the same figure in the old units) is cleared out entirely.
$ws.Range("Z23").Value  = 2.4512326436167768
$ws.Range("AA23").Value = 2.0191223235726747
$ws.Range("AB23").ClearContents()

$ws.Range("Z24").Value  = 2.7289748820247333
$ws.Range("AA24").Value = 1.3336153655648415
$ws.Range("AB24").ClearContents()

$ws.Range("Z25").Value  = 2.9308764540574193
$ws.Range("AA25").Value = 0.71375630065532591
$ws.Range("AB25").ClearContents()

$ws.Range("Z26").Value  = 3.2920359298333337
$ws.Range("AA26").Value = 0.40944545194439969
$ws.Range("AB26").ClearContents()

# New empty marker cell AD10, carrying the same "#,##0" number format
# already used for the Y23:Y26 id cells (cellXfs index 3).
$ws.Range("AD10").NumberFormat = "#,##0"

# Leave the cursor where the author left it on save.
$ws.Range("Y30").Select()
